$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 14
$ws.Range("A14").NumberFormat = "@"
$ws.Range("A14").Value = "2023-07-24"
$ws.Range("A14").ClearFormats()

$ws.Range("B14").Value = "RENGA STORE"
$ws.Range("C14").Value = 15011

# Row 15
$ws.Range("A15").NumberFormat = "@"
$ws.Range("A15").Value = "2023-07-24"
$ws.Range("A15").ClearFormats()

$ws.Range("B15").Value = "THIRUPATHI STORE"
$ws.Range("C15").Value = 23791
